$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99, shifting existing rows 99:161 down to 100:162.
$ws.Rows("99:99").Insert()

# Populate the newly inserted row 99 with the new data record.
$ws.Range("A99").Value = 11
$ws.Range("B99").Value = "Vega Monumental Concepción"
$ws.Range("C99").Value = "Bíobío"
$ws.Range("D99").Value = 45086
$ws.Range("E99").Value = 8
$ws.Range("F99").Value = 100112001
$ws.Range("G99").Value = "Berenjena"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 170
$ws.Range("K99").Value = 8000
$ws.Range("L99").Value = 9000
$ws.Range("M99").Value = 8529
$ws.Range("N99").Value = '$/caja 50 unidades'
$ws.Range("O99").Value = "Región de Arica y Parinacota"
$ws.Range("P99").Value = 171
$ws.Range("Q99").Value = 50
$ws.Range("R99").Value = "Hortaliza"

# Match the date-format style used by the rest of column D.
$ws.Range("D99").NumberFormat = $ws.Range("D100").NumberFormat
